$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capacitor C3 (10uF, row 55) was dropped from the design; C3 has effectively
# merged with C2 (both "Cap C" / 4.7 uF) per the firmware/board fix described
# in the commit. Delete the old "10uF"/"C3" row - this shifts the "Cap C"/"C2"
# row up from 56 to 55, and the "220uF 50v"/"C1, C4" row up from 57 to 56.
$ws.Rows.Item(55).Delete()

# The row that shifted into 55 ("Cap C" / "C2" / "4.7 uF" / qty 1) now also
# covers C3, so its designator list and quantity need updating. The leading
# apostrophe preserves the text/quote-prefix formatting the designator column
# already used (matches style of the surrounding designator cells).
$ws.Range("B55").Value = "'C2, C3"
$ws.Range("D55").Value = 2

$wb.Save()
